$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("M9").Value = 7000
$ws.Range("N9").Value = 2468.38
$ws.Range("O9").Value = 2468.38

# Row 10
$ws.Range("N10").Value = 4000.98

# Row 12
$ws.Range("N12").Value = 380804.16

# Row 19
$ws.Range("K19").Value = 10553.33
